# SSDM-12286 Added support for packaging entity validation plugins and dynamic property plugins.
# Replace the inline Python validation script with a reference to an external file "test.py".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "test.py"

# Update the active selection to C3 (as left by the edit in the source workbook).
$ws.Range("C3").Select()
